$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ B = 0.003078177322033415; C = 0.002658071450198252; D = 18.71679738969934;  E = 13.86384647080068;  G = 32.58638010927226 }
    3 = @{ B = 0.003078177322033415; C = 0.002658071450198252; D = 0.1496068669990043; E = 0.5333859586016987;  G = 0.6887290743729346 }
    4 = @{ B = 0.2881169905109251;   C = 0.3048912486333797;   D = 189.6080260415259;  E = 13.86384647080068;  G = 204.0648807514709 }
    5 = @{ B = 0.6545652718822623;   C = 1.626987699542094;    D = 0.7210945179870265; E = 0.5333859586016987;  G = 3.536033448013082 }
    6 = @{ B = 3.272327238179451;    C = 1.626987699542094;    D = 0.7210945179870265; E = 0.5333859586016987;  G = 6.15379541431027 }
    7 = @{ B = 0.04172184405617529;  C = 0.04103571897497393;  D = 0.7210945179870265; E = 0.5333859586016987;  G = 1.337238039619874 }
    8 = @{ B = 3.272327238179451;    C = 1.626987699542094;    D = 0.1496068669990043; E = 0.5333859586016987;  G = 5.582307763322248 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("G$row").Value = $vals.G
}
